# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Single apostrophe used as a text-prefix so numeric-looking strings
# (e.g. "1.00", "12.20") are stored as text, not auto-converted numbers.
$apo = "'"

$ws.Range('D2').Value = '67.327.99'
$ws.Range('D3').Value = '3.718.27'
$ws.Range('E3').Value = '  -3.28%  '
$ws.Range('D4').Value = $apo + '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').Value = $apo + '597.69'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.80%  '
$ws.Range('D6').Value = $apo + '166.58'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.66%  '
$ws.Range('D7').Value = '3.716.75'
$ws.Range('E7').Value = '  -3.46%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = $apo + '0.533'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.27%  '
$ws.Range('E10').Value = '  +1.20%  '
$ws.Range('D11').Value = $apo + '6.19'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.84%  '
$ws.Range('E12').Value = '  -3.34%  '
$ws.Range('D13').Value = $apo + '37.76'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.99%  '
$ws.Range('D14').Value = $apo + '0.0000244'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.74%  '
$ws.Range('D15').Value = '4.336.54'
$ws.Range('E15').Value = '  -3.27%  '
$ws.Range('D16').Value = '3.721.14'
$ws.Range('E16').Value = '  -3.01%  '
$ws.Range('D17').Value = '67.380.20'
$ws.Range('E17').Value = '  -2.92%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').Value = $apo + '7.27'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.20%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = $apo + '17.65'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +7.61%  '
$ws.Range('D20').Value = $apo + '0.115'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.14%  '
$ws.Range('D21').Value = $apo + '487.45'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.68%  '
$ws.Range('E22').Value = '  -3.20%  '
$ws.Range('D23').Value = $apo + '0.728'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.29%  '
$ws.Range('D24').Value = $apo + '85.11'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.33%  '
$ws.Range('D25').Value = $apo + '0.0000142'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.46%  '
$ws.Range('E26').Value = '  -4.68%  '
$ws.Range('D27').Value = $apo + '12.20'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.81%  '
$ws.Range('E28').Value = '  -2.28%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('E30').Value = '  -1.48%  '
$ws.Range('E31').Value = '  -6.52%  '
$ws.Range('D32').Value = $apo + '7.69'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.38%  '
$ws.Range('D33').Value = $apo + '31.50'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.59%  '
$ws.Range('D34').Value = '3.855.73'
$ws.Range('E34').Value = '  -3.37%  '
$ws.Range('E35').Value = '  -4.18%  '
$ws.Range('D36').Value = '3.658.60'
$ws.Range('E36').Value = '  -2.97%  '
$ws.Range('D37').Value = $apo + '0.999'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.21%  '
$ws.Range('D38').Value = $apo + '0.998'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.38%  '
$ws.Range('E39').Value = '  -3.37%  '
$ws.Range('E40').Value = '  -5.76%  '
$ws.Range('E41').Value = '  -2.80%  '
$ws.Range('D42').Value = $apo + '48.72'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.53%  '
$ws.Range('D43').Value = $apo + '428.39'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.66%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').Value = $apo + '2.83'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.81%  '
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').Value = $apo + '1.93'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.87%  '
$ws.Range('E46').Value = '  -0.60%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').Value = $apo + '40.46'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.50%  '
$ws.Range('D49').Value = $apo + '140.59'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.22%  '
$ws.Range('E50').Value = '  -2.37%  '
$ws.Range('D51').Value = '2.749.43'
$ws.Range('E51').Value = '  -4.94%  '